# Correct the arguments and fields for static tract data.
#
# Each of the 7 per-BlockGroup data sheets (BlkGrpecon, BlkGrpmedhhinc,
# BlkGrpeduc, BlkGrplatin, BlkGrppoverty, BlkGrpcommute, BlkGrprace) had a
# row 2 that documented a "GEOID / Id / GEO_ID" field which doesn't belong
# in the per-sheet tables (it's already covered by the STATE/COUNTY/TRACT/
# BlkGrp rows that directly follow it). Remove that row from every sheet,
# letting Excel shift the remaining rows up, and leave the selection on
# the new row 2 (matching what Excel leaves selected right after a row
# delete). Finally, land on the last data sheet (BlkGrprace).

$wb = $excel.ActiveWorkbook

$sheetNames = @(
    "BlkGrpecon",
    "BlkGrpmedhhinc",
    "BlkGrpeduc",
    "BlkGrplatin",
    "BlkGrppoverty",
    "BlkGrpcommute",
    "BlkGrprace"
)

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Rows(2).Delete()
    $ws.Activate()
    $ws.Range("A2:XFD2").Select()
}
